$d = $word.ActiveDocument

# The document contains several "<id>...</id>" sequences that were each
# split across three separate runs (one run for the literal "<id>" text,
# one for the identifier text, and one for the literal "</id>" text).
# This edit merges each trio of runs into a single run (keeping the
# formatting of the first "<id>" run) for the four ids that belong to
# the newly downloaded tc/tcn/tl content: p048r_1, p048r_2, p048r_3 and
# p048r_4. The unrelated "fig_p048r_1" id (used inside a <figure> block)
# is intentionally left untouched.

$ids = @("p048r_1", "p048r_2", "p048r_3", "p048r_4")

foreach ($id in $ids) {
    $searchText = "<id>" + $id + "</id>"
    $markerText = "@@MERGE_MARKER_" + $id + "@@"

    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

    if ($found) {
        # Assigning different text than what is currently in the range
        # forces Word to collapse the (three) runs the range spans into a
        # single run, inheriting the formatting of the first original run.
        # Re-assigning the exact original text back is treated as a no-op
        # by the engine and would NOT merge the runs, so we first swap in
        # a temporary marker string and then replace that marker with the
        # real text.
        $rng.Text = $markerText

        $rng2 = $d.Content
        $rng2.Find.ClearFormatting()
        $found2 = $rng2.Find.Execute($markerText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found2) {
            $rng2.Text = $searchText
        }
    }
}
